$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "65.862.45"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "3.582.63"
$ws.Range("E3").Value = "  +1.04%  "
Set-TextValue $ws "D4" "1.00"
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue $ws "D5" "601.04"
$ws.Range("E5").Value = "  +0.47%  "
Set-TextValue $ws "D6" "137.57"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").Value = "3.583.66"
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("E8").Value = "  -0.08%  "
Set-TextValue $ws "D9" "0.499"
$ws.Range("E9").Value = "  +0.79%  "
Set-TextValue $ws "D10" "0.125"
$ws.Range("E10").Value = "  +0.48%  "
Set-TextValue $ws "D11" "7.19"
$ws.Range("E11").Value = "  +3.82%  "
Set-TextValue $ws "D12" "0.391"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "4.205.05"
$ws.Range("E13").Value = "  +1.37%  "
Set-TextValue $ws "D14" "27.96"
$ws.Range("E14").Value = "  +2.35%  "
Set-TextValue $ws "D15" "0.0000184"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "3.594.18"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "64.846.25"
$ws.Range("E18").Value = "  -0.43%  "
Set-TextValue $ws "D19" "9.94"
$ws.Range("E19").Value = "  -1.25%  "
Set-TextValue $ws "D20" "14.59"
$ws.Range("E20").Value = "  +2.39%  "
Set-TextValue $ws "D21" "5.81"
$ws.Range("E21").Value = "  -1.20%  "
Set-TextValue $ws "D22" "395.46"
$ws.Range("E22").Value = "  +0.79%  "
Set-TextValue $ws "D23" "0.585"
$ws.Range("E23").Value = "  +1.91%  "
$ws.Range("D24").Value = "3.739.73"
$ws.Range("E24").Value = "  +1.44%  "
Set-TextValue $ws "D25" "74.98"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("E26").Value = "  +0.03%  "
Set-TextValue $ws "D27" "0.0000118"
$ws.Range("E27").Value = "  +3.27%  "
Set-TextValue $ws "D28" "8.01"
$ws.Range("E28").Value = "  +1.86%  "
Set-TextValue $ws "D29" "1.64"
$ws.Range("E29").Value = "  +19.40%  "
Set-TextValue $ws "D30" "8.58"
$ws.Range("E30").Value = "  +3.21%  "
Set-TextValue $ws "D31" "1.00"
$ws.Range("E31").Value = "  +0.36%  "
Set-TextValue $ws "D32" "2.32"
$ws.Range("E32").Value = "  +1.86%  "
$ws.Range("D33").Value = "3.601.10"
Set-TextValue $ws "D34" "24.47"
$ws.Range("E34").Value = "  +2.58%  "
Set-TextValue $ws "D35" "0.148"
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("E36").Value = "  +0.01%  "
Set-TextValue $ws "D37" "5.36"
$ws.Range("E37").Value = "  +6.54%  "
Set-TextValue $ws "D38" "1.60"
$ws.Range("E38").Value = "  +1.28%  "
Set-TextValue $ws "D39" "6.96"
$ws.Range("E39").Value = "  +0.35%  "
Set-TextValue $ws "D40" "169.81"
$ws.Range("E40").Value = "  +0.35%  "
Set-TextValue $ws "D41" "0.0834"
$ws.Range("E41").Value = "  +3.59%  "
Set-TextValue $ws "D42" "0.838"
$ws.Range("E42").Value = "  +1.75%  "
Set-TextValue $ws "D43" "26.26"
$ws.Range("E43").Value = "  -1.14%  "
Set-TextValue $ws "D44" "1.26"
$ws.Range("E44").Value = "  +4.69%  "
Set-TextValue $ws "D45" "43.05"
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws "D46" "1.00"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D47" "4.51"
$ws.Range("E47").Value = "  +1.39%  "
Set-TextValue $ws "D48" "1.69"
$ws.Range("E48").Value = "  +1.12%  "
Set-TextValue $ws "D49" "6.91"
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("D50").Value = "2.440.66"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws "D51" "0.905"
$ws.Range("E51").Value = "  +6.74%  "
